$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '26.196.47'
$ws.Range("E2").Value = '  -1.96%  '
$ws.Range("D3").Value = '1.581.14'
$ws.Range("E3").Value = '  -1.28%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").Value = '''209.37'
$ws.Range("E5").Value = '  -1.08%  '
$ws.Range("D6").Value = '''0.502'
$ws.Range("E6").Value = '  -2.10%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").Value = '''0.0610'
$ws.Range("E8").Value = '  -1.52%  '
$ws.Range("D9").Value = '''0.245'
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("D10").Value = '''19.51'
$ws.Range("E10").Value = '  -0.82%  '
$ws.Range("D11").Value = '''0.0846'
$ws.Range("E11").Value = '  +0.26%  '
$ws.Range("D12").Value = '1.803.07'
$ws.Range("E12").Value = '  -1.28%  '
$ws.Range("D13").Value = '1.598.86'
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("D14").Value = '''4.03'
$ws.Range("E14").Value = '  -0.16%  '
$ws.Range("E15").Value = '  -1.28%  '
$ws.Range("E16").Value = '  -0.55%  '
$ws.Range("D17").Value = '26.190.74'
$ws.Range("E17").Value = '  -1.87%  '
$ws.Range("E18").Value = '  -1.24%  '
$ws.Range("D19").Value = '''7.25'
$ws.Range("E19").Value = '  +0.94%  '
$ws.Range("E20").Value = '  -0.22%  '
$ws.Range("D21").Value = '''206.35'
$ws.Range("E21").Value = '  -1.83%  '
$ws.Range("D22").Value = '''4.26'
$ws.Range("E22").Value = '  -0.59%  '
$ws.Range("E23").Value = '  -3.34%  '
$ws.Range("E24").Value = '  -1.37%  '
$ws.Range("D25").Value = '''144.92'
$ws.Range("E25").Value = '  +0.44%  '
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("D27").Value = '''6.99'
$ws.Range("E27").Value = '  -1.59%  '
$ws.Range("E28").Value = '  -1.08%  '
$ws.Range("E29").Value = '  -1.08%  '
$ws.Range("E30").Value = '  -1.97%  '
$ws.Range("E31").Value = '  -1.22%  '
$ws.Range("D32").Value = '''3.21'
$ws.Range("E32").Value = '  -1.48%  '
$ws.Range("E33").Value = '  -1.23%  '
$ws.Range("D34").Value = '''1.26'
$ws.Range("E34").Value = '  +6.59%  '
$ws.Range("D35").Value = '1.282.56'
$ws.Range("E35").Value = '  -1.22%  '
$ws.Range("E36").Value = '  -0.22%  '
$ws.Range("D37").Value = '''0.604'
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("E38").Value = '  -1.14%  '
$ws.Range("D39").Value = '''0.0166'
$ws.Range("E39").Value = '  -1.84%  '
$ws.Range("E40").Value = '  -1.74%  '
$ws.Range("D41").Value = '''5.51'
$ws.Range("E41").Value = '  +2.01%  '
$ws.Range("E42").Value = '  -0.97%  '
$ws.Range("D43").Value = '''62.44'
$ws.Range("E43").Value = '  -0.94%  '
$ws.Range("E44").Value = '  -3.08%  '
$ws.Range("D45").Value = '1.716.18'
$ws.Range("E45").Value = '  -1.41%  '
$ws.Range("D46").Value = '''88.62'
$ws.Range("E46").Value = '  -2.18%  '
$ws.Range("E47").Value = '  -0.60%  '
$ws.Range("E48").Value = '  -0.08%  '
$ws.Range("E49").Value = '  -1.65%  '
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("D51").Value = '''7.37'
$ws.Range("E51").Value = '  -0.87%  '

Write-Output "Updated 77 cells across D/E columns for cryptos sheet"
